{"js": "// Change 1: \"Batas Akhir Penawaran\" line \u2014 insert \"WIB \" right before\n// \"waktu server lelang)\" so it reads \"... WIB waktu server lelang) \".\nconst target1 = context.document.body.search(\"waktu server lelang) \", { matchCase: true });\ntarget1.load(\"items\");\nawait context.sync();\n\nif (target1.items.length === 0) {\n  throw new Error('Target text \"waktu server lelang) \" not found');\n}\ntarget1.items[0].insertText(\"WIB \", Word.InsertLocation.before);\nawait context.sync();\n\n// Change 2: remove the leading \"sebagai pengumuman lelang\" phrase so the\n// sentence goes straight from \"${tanggalPengumuman}\" to the comma.\nconst target2 = context.document.body.search(\" sebagai pengumuman lelang\", { matchCase: true });\ntarget2.load(\"items\");\nawait context.sync();\n\nif (target2.items.length === 0) {\n  throw new Error('Target text \" sebagai pengumuman lelang\" not found');\n}\ntarget2.items[0].delete();\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n#   MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n#   ReplaceWith, Replace)\n\n# Change 1: \"Batas Akhir Penawaran\" line \u2014 insert \"WIB \" right before\n# \"waktu server lelang)\" so it reads \"... WIB waktu server lelang) \".\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$found1 = $rng1.Find.Execute(\"waktu server lelang) \", $false, $false, $false, $false, $false, $true, $true, $false, \"\", $false)\nif (-not $found1) {\n    throw 'Target text \"waktu server lelang) \" not found'\n}\n$rng1.InsertBefore(\"WIB \")\n\n# Change 2: remove the leading \"sebagai pengumuman lelang\" phrase so the\n# sentence goes straight from \"${tanggalPengumuman}\" to the comma.\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$found2 = $rng2.Find.Execute(\" sebagai pengumuman lelang\", $false, $false, $false, $false, $false, $true, $true, $false, \"\", $false)\nif (-not $found2) {\n    throw 'Target text \" sebagai pengumuman lelang\" not found'\n}\n$rng2.Delete()\n"}
